$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")

# Row 40
$ws.Range("H40").Value = 2414.7693
$ws.Range("I40").Value = 4767
$ws.Range("J40").Value = 1709.1
$ws.Range("K40").Value = 4767
$ws.Range("L40").Value = 1709.1
$ws.Range("M40").Value = -4592
$ws.Range("N40").Value = -2059.1

# Row 106
$ws.Range("H106").Value = 9596.066000000001
$ws.Range("I106").Value = 11111.25
$ws.Range("K106").Value = 11111.25
$ws.Range("M106").Value = -10480.25

# Row 128
$ws.Range("H128").Value = 29666
$ws.Range("J128").Value = 29666
$ws.Range("L128").Value = 29666
$ws.Range("N128").Value = -39626

# Row 132
$ws.Range("H132").Value = 3902
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# Row 137
$ws.Range("H137").Value = 1237.2759
$ws.Range("I137").Value = 1043.4706
$ws.Range("J137").Value = 1511.8334
$ws.Range("K137").Value = 3130.4118
$ws.Range("L137").Value = 4535.5002
$ws.Range("M137").Value = -580.4118000000003
$ws.Range("N137").Value = -9635.5002


# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 2995.3225
$ws.Range("I32").Value = 2959.926
$ws.Range("J32").Value = 3234.25
$ws.Range("K32").Value = 2959.926
$ws.Range("L32").Value = 3234.25
$ws.Range("M32").Value = -2672.926
$ws.Range("N32").Value = -3808.25

# Row 45
$ws.Range("H45").Value = 1468.3572
$ws.Range("I45").Value = 1596.091
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 1596.091
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -1219.091
$ws.Range("N45").Value = -1754

# Row 61
$ws.Range("H61").Value = 1523.2
$ws.Range("I61").Value = 1205.5555
$ws.Range("J61").Value = 1999.6666
$ws.Range("K61").Value = 1205.5555
$ws.Range("L61").Value = 1999.6666
$ws.Range("M61").Value = -993.5554999999999
$ws.Range("N61").Value = -2423.6666

# Row 74
$ws.Range("H74").Value = 1619.0435
$ws.Range("I74").Value = 674.9167
$ws.Range("J74").Value = 2649
$ws.Range("K74").Value = 674.9167
$ws.Range("L74").Value = 2649
$ws.Range("M74").Value = 199.0833
$ws.Range("N74").Value = -4397

# Row 77
$ws.Range("H77").Value = 1619.0435
$ws.Range("I77").Value = 674.9167
$ws.Range("J77").Value = 2649
$ws.Range("K77").Value = 3374.5835
$ws.Range("L77").Value = 13245
$ws.Range("M77").Value = 993.4165000000003
$ws.Range("N77").Value = -21981

# Row 122
$ws.Range("H122").Value = 1138.4546
$ws.Range("I122").Value = 820.375
$ws.Range("K122").Value = 2461.125
$ws.Range("M122").Value = -11.125

# Row 132
$ws.Range("H132").Value = 1976.6842
$ws.Range("I132").Value = 1686.7333
$ws.Range("K132").Value = 5060.199900000001
$ws.Range("M132").Value = -2530.199900000001

# Row 136
$ws.Range("H136").Value = 1523.2
$ws.Range("I136").Value = 1205.5555
$ws.Range("J136").Value = 1999.6666
$ws.Range("K136").Value = 3616.6665
$ws.Range("L136").Value = 5998.9998
$ws.Range("M136").Value = -1066.6665
$ws.Range("N136").Value = -11098.9998


# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")

# Row 134
$ws.Range("H134").Value = 7489.05
$ws.Range("I134").Value = 1252.2
$ws.Range("J134").Value = 26199.6
$ws.Range("K134").Value = 3756.6
$ws.Range("L134").Value = 78598.79999999999
$ws.Range("M134").Value = -1221.6
$ws.Range("N134").Value = -83668.79999999999


# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")

# Row 18
$ws.Range("H18").Value = 47795
$ws.Range("J18").Value = 47795
$ws.Range("L18").Value = 47795
$ws.Range("N18").Value = -48255

# Row 31
$ws.Range("H31").Value = 1496.6666
$ws.Range("I31").Value = 1075.5883
$ws.Range("J31").Value = 2212.5
$ws.Range("K31").Value = 1075.5883
$ws.Range("L31").Value = 2212.5
$ws.Range("M31").Value = -780.5882999999999
$ws.Range("N31").Value = -2802.5

# Row 34
$ws.Range("H34").Value = 1496.6666
$ws.Range("I34").Value = 1075.5883
$ws.Range("J34").Value = 2212.5
$ws.Range("K34").Value = 1075.5883
$ws.Range("L34").Value = 2212.5
$ws.Range("M34").Value = -873.5882999999999
$ws.Range("N34").Value = -2616.5

# Row 58
$ws.Range("H58").Value = 1434.6875
$ws.Range("I58").Value = 1244.6
$ws.Range("K58").Value = 1244.6
$ws.Range("M58").Value = -1041.6

# Row 129
$ws.Range("H129").Value = 49499.25
$ws.Range("J129").Value = 49499.25
$ws.Range("L129").Value = 49499.25
$ws.Range("N129").Value = -59499.25

# Row 132
$ws.Range("H132").Value = 5284.7417
$ws.Range("I132").Value = 6599.4
$ws.Range("K132").Value = 19798.2
$ws.Range("M132").Value = -17268.2

# Row 134
$ws.Range("H134").Value = 2229.3914
$ws.Range("I134").Value = 2461.4666
$ws.Range("J134").Value = 1794.25
$ws.Range("K134").Value = 7384.399800000001
$ws.Range("L134").Value = 5382.75
$ws.Range("M134").Value = -4849.399800000001
$ws.Range("N134").Value = -10452.75

# Row 136
$ws.Range("H136").Value = 1434.6875
$ws.Range("I136").Value = 1244.6
$ws.Range("K136").Value = 3733.8
$ws.Range("M136").Value = -1183.8


# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")

# Row 34
$ws.Range("H34").Value = 30000
$ws.Range("J34").Value = 30000
$ws.Range("L34").Value = 30000
$ws.Range("N34").Value = -30536

# Row 70
$ws.Range("H70").Value = 37503556
$ws.Range("I70").Value = 31253712
$ws.Range("J70").Value = 50003250
$ws.Range("K70").Value = 31253712
$ws.Range("L70").Value = 50003250
$ws.Range("M70").Value = -31253442
$ws.Range("N70").Value = -50003790

# Row 73
$ws.Range("H73").Value = 37503556
$ws.Range("I73").Value = 31253712
$ws.Range("J73").Value = 50003250
$ws.Range("K73").Value = 31253712
$ws.Range("L73").Value = 50003250
$ws.Range("M73").Value = -31252776
$ws.Range("N73").Value = -50005122

# Row 76
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630

# Row 79
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184

# Row 122
$ws.Range("H122").Value = 1451.8518
$ws.Range("I122").Value = 1584.65
$ws.Range("J122").Value = 1072.4286
$ws.Range("K122").Value = 4753.950000000001
$ws.Range("L122").Value = 3217.2858
$ws.Range("M122").Value = -2303.950000000001
$ws.Range("N122").Value = -8117.2858

# Row 123
$ws.Range("H123").Value = 10325.909
$ws.Range("J123").Value = 10325.909
$ws.Range("L123").Value = 10325.909
$ws.Range("N123").Value = -15225.909

# Row 132
$ws.Range("H132").Value = 2083.2666
$ws.Range("I132").Value = 1919.1904
$ws.Range("J132").Value = 2466.111
$ws.Range("K132").Value = 5757.5712
$ws.Range("L132").Value = 7398.333
$ws.Range("M132").Value = -3227.5712
$ws.Range("N132").Value = -12458.333


# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 976.8125
$ws.Range("I22").Value = 518.125
$ws.Range("J22").Value = 1435.5
$ws.Range("K22").Value = 518.125
$ws.Range("L22").Value = 1435.5
$ws.Range("M22").Value = -223.125
$ws.Range("N22").Value = -2025.5

# Row 27
$ws.Range("H27").Value = 976.8125
$ws.Range("I27").Value = 518.125
$ws.Range("J27").Value = 1435.5
$ws.Range("K27").Value = 518.125
$ws.Range("L27").Value = 1435.5
$ws.Range("M27").Value = -411.125
$ws.Range("N27").Value = -1649.5

# Row 82
$ws.Range("H82").Value = 1216.8334
$ws.Range("I82").Value = 1100.2222
$ws.Range("K82").Value = 1100.2222
$ws.Range("M82").Value = -739.2221999999999

# Row 85
$ws.Range("H85").Value = 1216.8334
$ws.Range("I85").Value = 1100.2222
$ws.Range("K85").Value = 1100.2222
$ws.Range("M85").Value = 147.7778000000001

# Row 122
$ws.Range("H122").Value = 13164774
$ws.Range("I122").Value = 17865734
$ws.Range("K122").Value = 53597202
$ws.Range("M122").Value = -53594752

# Row 132
$ws.Range("H132").Value = 18208.717
$ws.Range("I132").Value = 1175.4474
$ws.Range("J132").Value = 47629.816
$ws.Range("K132").Value = 3526.3422
$ws.Range("L132").Value = 142889.448
$ws.Range("M132").Value = -996.3422
$ws.Range("N132").Value = -147949.448

# Row 136
$ws.Range("H136").Value = 4599
$ws.Range("I136").Value = 5094.6523
$ws.Range("J136").Value = 799
$ws.Range("K136").Value = 15283.9569
$ws.Range("L136").Value = 2397
$ws.Range("M136").Value = -12733.9569
$ws.Range("N136").Value = -7497


# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")

# Row 52
$ws.Range("H52").Value = 15399.667
$ws.Range("J52").Value = 16599.5
$ws.Range("L52").Value = 16599.5
$ws.Range("N52").Value = -17051.5

# Row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

# Row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

# Row 119
$ws.Range("H119").Value = 18978.8
$ws.Range("J119").Value = 18978.8
$ws.Range("L119").Value = 18978.8
$ws.Range("N119").Value = -28654.8

# Row 122
$ws.Range("H122").Value = 23638532
$ws.Range("I122").Value = 23638532
$ws.Range("K122").Value = 70915596
$ws.Range("M122").Value = -70913146

# Row 132
$ws.Range("H132").Value = 2625.72
$ws.Range("I132").Value = 2607.9092
$ws.Range("K132").Value = 7823.7276
$ws.Range("M132").Value = -5293.7276

